$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 2.05
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 2.88
$ws.Range("L3").Value = 4.75
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.85
$ws.Range("U3").Value = 4.3
$ws.Range("V3").Value = 1.21
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 1.14
$ws.Range("AD3").Value = 8.5
$ws.Range("AE3").Value = 10
$ws.Range("AF3").Value = 19
$ws.Range("AK3").Value = 19
$ws.Range("AN3").Value = 8.5
$ws.Range("AO3").Value = 19
$ws.Range("AP3").Value = 15
$ws.Range("AQ3").Value = 41
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("U4").Value = 3.35
$ws.Range("V4").Value = 1.31
$ws.Range("G5").Value = 2.9
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 2.15
$ws.Range("J5").Value = 3.35
$ws.Range("K5").Value = 2.22
$ws.Range("L5").Value = 2.7
$ws.Range("O5").Value = 1.21
$ws.Range("P5").Value = 3.5
$ws.Range("S5").Value = 1.65
$ws.Range("T5").Value = 1.98
$ws.Range("W5").Value = 2.52
$ws.Range("X5").Value = 1.4
$ws.Range("AA5").Value = 1.57
$ws.Range("AB5").Value = 2.1
$ws.Range("AC5").Value = 11.25
$ws.Range("AD5").Value = 16.5
$ws.Range("AE5").Value = 10.5
$ws.Range("AF5").Value = 35
$ws.Range("AG5").Value = 22
$ws.Range("AH5").Value = 27
$ws.Range("AI5").Value = 13
$ws.Range("AJ5").Value = 7.1
$ws.Range("AK5").Value = 13
$ws.Range("AM5").Value = 300
$ws.Range("AN5").Value = 9.25
$ws.Range("AO5").Value = 11.5
$ws.Range("AP5").Value = 8.75
$ws.Range("AQ5").Value = 21
$ws.Range("AR5").Value = 16
$ws.Range("AS5").Value = 23
